$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: refresh the rolled-up Status column (B & C) for both
# rows now that handback is complete.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# Helper: mirror the source-file / target-file links into the new
# "Latest Target File" (F) and "Latest Handback File" (G) columns, and
# stamp the handback datetime (H), for one language sheet.
# ---------------------------------------------------------------------
function Update-LangSheet([string]$sheetName, [string]$xlfFile, [string]$xlfUrl, [string]$handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    $mdFile = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.md"
    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/eb727aa482a48db61efa329cd9e4c530cdd4fdcc/e2e/$mdFile"

    # Row 2
    $ws.Range("C2").Value = $statusText
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl, "", "", $mdFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl, "", "", $xlfFile) | Out-Null
    $ws.Range("F2").Style = "HyperLink"
    $ws.Range("G2").Style = "HyperLink"
    $ws.Range("H2").Value = $handbackDateTime

    # Row 3
    $ws.Range("C3").Value = $statusText
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrl, "", "", $mdFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrl, "", "", $xlfFile) | Out-Null
    $ws.Range("F3").Style = "HyperLink"
    $ws.Range("G3").Style = "HyperLink"
    $ws.Range("H3").Value = $handbackDateTime
}

$zhCnXlfFile = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.zh-cn.xlf"
$zhCnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16301a506697b8e78e76f61f92c639f869e67655/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhCnXlfFile"

$deDeXlfFile = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.de-de.xlf"
$deDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4ce29a537461e301605659247a25f3005f538009/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deDeXlfFile"

Update-LangSheet "zh-cn" $zhCnXlfFile $zhCnXlfUrl "2016-03-24 03:14:59"
Update-LangSheet "de-de" $deDeXlfFile $deDeXlfUrl "2016-03-24 03:15:08"
